$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = 1814
    "E2" = 307
    "F2" = 307
    "G2" = 367
    "H2" = 133
    "I2" = 95
    "J2" = 38
    "K2" = 3027
    "L2" = 1142
    "M2" = 1885
    "N2" = 1580
    "O2" = 305
    "P2" = 91
    "Q2" = 36
    "R2" = -165
    "S2" = 33
    "T2" = 141
    "U2" = -105
    "V2" = 473
    "W2" = 16.94
    "X2" = 7.33
    "Y2" = 6.17
    "Z2" = 4.59
    "AA2" = 60.59
    "AB2" = 1590.8
    "AC2" = 2598
    "AD2" = 23.01
    "AE2" = 47078
    "AF2" = 1.27
    "AG2" = 700
    "AH2" = 1.17
    "AI2" = 24.83
    "AJ2" = 3640950
    "D3" = 2218
    "E3" = 282
    "F3" = 282
    "G3" = 330
    "H3" = 306
    "I3" = 152
    "J3" = 154
    "K3" = 5010
    "L3" = 1198
    "M3" = 3812
    "N3" = 2929
    "O3" = 884
    "P3" = 112
    "Q3" = 82
    "R3" = -498
    "S3" = 862
    "T3" = 285
    "U3" = -203
    "V3" = 493
    "W3" = 12.7
    "X3" = 13.78
    "Y3" = 6.75
    "Z3" = 7.61
    "AA3" = 31.41
    "AB3" = 2429.46
    "AC3" = 3644
    "AD3" = 30.46
    "AE3" = 69844
    "AF3" = 1.59
    "AG3" = 1200
    "AH3" = 1.08
    "AI3" = 33.09
    "AJ3" = 4477766
    "D4" = 3357
    "E4" = 230
    "F4" = 230
    "G4" = 489
    "H4" = 429
    "I4" = 294
    "J4" = 134
    "K4" = 7057
    "L4" = 1576
    "M4" = 5480
    "N4" = 3604
    "O4" = 1877
    "P4" = 125
    "Q4" = 194
    "R4" = -123
    "S4" = -51
    "T4" = 265
    "U4" = -71
    "V4" = 621
    "W4" = 6.87
    "X4" = 12.77
    "Y4" = 9.01
    "Z4" = 7.1
    "AA4" = 28.77
    "AB4" = 2797.22
    "AC4" = 6225
    "AD4" = 11.07
    "AE4" = 76263
    "AF4" = 0.9
    "AG4" = 1400
    "AH4" = 2.03
    "AI4" = 22.48
    "AJ4" = 5009861
    "D5" = 4507
    "E5" = 435
    "F5" = 435
    "G5" = 792
    "H5" = 674
    "I5" = 480
    "J5" = 194
    "K5" = 8189
    "L5" = 1787
    "M5" = 6402
    "N5" = 4089
    "O5" = 2313
    "P5" = 125
    "Q5" = 578
    "R5" = -487
    "S5" = -323
    "T5" = 163
    "U5" = 415
    "V5" = 474
    "W5" = 9.65
    "X5" = 14.96
    "Y5" = 12.48
    "Z5" = 8.84
    "AA5" = 27.91
    "AB5" = 3130.59
    "AC5" = 9583
    "AD5" = 7.57
    "AE5" = 81622
    "AF5" = 0.89
    "AG5" = 900
    "AH5" = 1.24
    "AI5" = 9.390000000000001
    "AJ5" = 5009861
    "D6" = 5311
    "E6" = 537
    "F6" = 537
    "G6" = 674
    "H6" = 528
    "I6" = 275
    "K6" = 8714
    "L6" = 2007
    "M6" = 6708
    "N6" = 4293
    "P6" = 125
    "Q6" = 379
    "R6" = -641
    "S6" = 53
    "T6" = 510
    "U6" = -131
    "V6" = 694
    "W6" = 10.11
    "X6" = 9.94
    "Y6" = 6.55
    "Z6" = 6.25
    "AA6" = 29.91
    "AB6" = 3313.49
    "AC6" = 5481
    "AD6" = 11.4
    "AE6" = 85685
    "AF6" = 0.73
    "AG6" = 1300
    "AH6" = 2.08
    "AI6" = 23.72
    "AJ6" = 5009861
    "D7" = 8793
    "E7" = 1792
    "G7" = 1907
    "H7" = 1526
    "I7" = 992
    "K7" = 11026
    "L7" = 2867
    "M7" = 8159
    "N7" = 5210
    "P7" = 125
    "Q7" = 1882
    "R7" = -1770
    "S7" = -66
    "T7" = 430
    "U7" = 1074
    "W7" = 20.38
    "X7" = 17.36
    "Y7" = 20.88
    "Z7" = 15.46
    "AA7" = 35.14
    "AC7" = 19801
    "AD7" = 5.68
    "AE7" = 103995
    "AF7" = 1.08
    "AG7" = 1500
    "AH7" = 1.33
    "AI7" = 7.58
    "D8" = 10109
    "E8" = 2247
    "G8" = 2393
    "H8" = 1914
    "I8" = 1283
    "K8" = 13182
    "L8" = 3192
    "M8" = 9990
    "N8" = 6409
    "P8" = 125
    "Q8" = 2153
    "R8" = -1046
    "S8" = -75
    "T8" = 512
    "U8" = 1395
    "W8" = 22.23
    "X8" = 18.93
    "Y8" = 22.09
    "Z8" = 15.81
    "AA8" = 31.95
    "AC8" = 25609
    "AD8" = 4.39
    "AE8" = 127928
    "AF8" = 0.88
    "AG8" = 1650
    "AH8" = 1.47
    "AI8" = 6.44
    "D9" = 11154
    "E9" = 2401
    "G9" = 2644
    "H9" = 2115
    "I9" = 1438
    "K9" = 15468
    "L9" = 3450
    "M9" = 12018
    "N9" = 7760
    "P9" = 125
    "Q9" = 2386
    "R9" = -670
    "S9" = -83
    "T9" = 214
    "U9" = 1840
    "W9" = 21.53
    "X9" = 18.96
    "Y9" = 20.3
    "Z9" = 14.76
    "AA9" = 28.71
    "AC9" = 28703
    "AD9" = 3.92
    "AE9" = 154895
    "AF9" = 0.73
    "AG9" = 1750
    "AH9" = 6.1
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

Write-Output "Updated $($updates.Count) cells"